$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$conv = $wb.Worksheets.Item("CONVERTION")

# Enter the number of days used for this pay period's conversion lookup on
# the CONVERTION helper sheet. This cascades (via existing formulas) into
# J4, K3, and L3 automatically.
$conv.Range("J3").Value = 15

# Record the new VL(10-0-0) earned-leave entry for the period starting
# 8/15/2023 (previously 8/1/2023) and its EARNED value.
$ws1.Range("A80").Value = 45153
$ws1.Range("C80").Value = 0.667

# Clear out the now-unused future period dates beyond row 80.
$ws1.Range("A81:A131").ClearContents()
